$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing date cell format (style index reused, numFmtId 14 short date)
# from C4 onto the C8:C10 range before setting their values.
$ws.Range("C4").Copy()
$ws.Range("C8:C10").PasteSpecial(-4122)  # xlPasteFormats

# Row 8
$ws.Range("A8").Value = "Анализ, формированеи рынков для организации ЛИГА-7"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 43512

# Row 9
$ws.Range("A9").Value = "Работа по созданию функционала загрузки данных"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 43510

# Row 10
$ws.Range("A10").Value = "Работа по созданию функционала загрузки данных"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 43512

# Update the sheet view selection to A10
$ws.Range("A10").Select()
